# Generate Report for Handback
# Update the handoff/handback timestamp strings recorded in the
# localization-status report for the latest handback generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D2").Value = "2016-52-17 16:52:45"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 16:52:41"
$wsZhCn.Range("H2").Value = "2016-03-17 16:52:59"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 16:52:45"
$wsDeDe.Range("H2").Value = "2016-03-17 16:53:10"
